$d = $word.ActiveDocument

$r = $d.Content
$r.Find.Execute("Frontend and Backend.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)

$full = " Then I deployed the exercise website to CSC Rathi and finished the tutorial."
$r.Text = $full

$sub = $d.Content
$sub.Find.Execute("the exercise website", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sub.Font.Kerning = 0
